$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13285
$ws1.Range("F6").Value = 436
$ws1.Range("F7").Value = 1299
$ws1.Range("F8").Value = 124

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13285
$ws4.Range("F8").Value = 436
$ws4.Range("F9").Value = 1299
$ws4.Range("F11").Value = 124
